$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59; existing rows 59:92 shift down to 60:93
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly data point
$ws.Cells.Item(59, 1).Value = 5
$ws.Cells.Item(59, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(59, 3).Value = "Maule"
$ws.Cells.Item(59, 4).Value = 44574
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = 100112030
$ws.Cells.Item(59, 7).Value = "Poroto granado"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 300
$ws.Cells.Item(59, 11).Value = 22000
$ws.Cells.Item(59, 12).Value = 22000
$ws.Cells.Item(59, 13).Value = 22000
$ws.Cells.Item(59, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(59, 15).Value = "Región del Maule"
$ws.Cells.Item(59, 16).Value = 880
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
